$d = $word.ActiveDocument

# Anchor paragraph that stays right before the block we want to remove.
$anchor = $d.Content
$anchor.Find.Execute("LOM3057: Introdução aos Materiais Poliméricos (Requisito fraco)") | Out-Null
$anchorPara = $anchor.Paragraphs(1)

# Anchor paragraph that is the last one we want to remove (the copyright line).
$lastToRemove = $d.Content
$lastToRemove.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution") | Out-Null
$lastToRemovePara = $lastToRemove.Paragraphs(1)

# Remove everything from right after the anchor paragraph's mark through the
# end (including the paragraph mark) of the copyright paragraph. That deletes
# the blank paragraph, the "Ver no Jupiter..." paragraph and the copyright
# paragraph in one pass, leaving the anchor paragraph intact and reconnecting
# it directly to whatever paragraph used to follow the copyright line.
$deleteRange = $d.Range($anchorPara.Range.End, $lastToRemovePara.Range.End)
$deleteRange.Delete()
